$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2122.4707
$ws.Range("I40").Value = 2075.6667
$ws.Range("K40").Value = 2075.6667
$ws.Range("M40").Value = -1900.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4594.7085
$ws.Range("I64").Value = 3620.7693
$ws.Range("J64").Value = 5745.727
$ws.Range("K64").Value = 3620.7693
$ws.Range("L64").Value = 5745.727
$ws.Range("M64").Value = -3372.7693
$ws.Range("N64").Value = -6241.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4594.7085
$ws.Range("I67").Value = 3620.7693
$ws.Range("J67").Value = 5745.727
$ws.Range("K67").Value = 3620.7693
$ws.Range("L67").Value = 5745.727
$ws.Range("M67").Value = -2762.7693
$ws.Range("N67").Value = -7461.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3367.5334
$ws.Range("I74").Value = 3250.6428
$ws.Range("J74").Value = 5004
$ws.Range("K74").Value = 3250.6428
$ws.Range("L74").Value = 5004
$ws.Range("M74").Value = -2314.6428
$ws.Range("N74").Value = -6876

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3271238.5
$ws.Range("I76").Value = 4118337.2
$ws.Range("J76").Value = 3857.1428
$ws.Range("K76").Value = 4118337.2
$ws.Range("L76").Value = 3857.1428
$ws.Range("M76").Value = -4118022.2
$ws.Range("N76").Value = -4487.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3367.5334
$ws.Range("I77").Value = 3250.6428
$ws.Range("J77").Value = 5004
$ws.Range("K77").Value = 16253.214
$ws.Range("L77").Value = 25020
$ws.Range("M77").Value = -11573.214
$ws.Range("N77").Value = -34380

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3271238.5
$ws.Range("I79").Value = 4118337.2
$ws.Range("J79").Value = 3857.1428
$ws.Range("K79").Value = 4118337.2
$ws.Range("L79").Value = 3857.1428
$ws.Range("M79").Value = -4117245.2
$ws.Range("N79").Value = -6041.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1049.7368
$ws.Range("J129").Value = 1091.4722
$ws.Range("L129").Value = 3274.4166
$ws.Range("N129").Value = -13274.4166

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 11822723
$ws.Range("I138").Value = 3969543
$ws.Range("J138").Value = 15154375
$ws.Range("K138").Value = 11908629
$ws.Range("L138").Value = 45463125
$ws.Range("M138").Value = -11903489
$ws.Range("N138").Value = -45473405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16650.291
$ws.Range("I32").Value = 4281.7075
$ws.Range("K32").Value = 4281.7075
$ws.Range("M32").Value = -3994.7075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1012
$ws.Range("I45").Value = 1012
$ws.Range("K45").Value = 1012
$ws.Range("M45").Value = -635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4215.9165
$ws.Range("I63").Value = 4438.8823
$ws.Range("J63").Value = 3674.4285
$ws.Range("K63").Value = 4438.8823
$ws.Range("L63").Value = 3674.4285
$ws.Range("M63").Value = -3752.8823
$ws.Range("N63").Value = -5046.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4215.9165
$ws.Range("I66").Value = 4438.8823
$ws.Range("J66").Value = 3674.4285
$ws.Range("K66").Value = 22194.4115
$ws.Range("L66").Value = 18372.1425
$ws.Range("M66").Value = -18762.4115
$ws.Range("N66").Value = -25236.1425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1000
$ws.Range("I88").Value = 1000
$ws.Range("K88").Value = 1000
$ws.Range("M88").Value = -594

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1000
$ws.Range("I91").Value = 1000
$ws.Range("K91").Value = 1000
$ws.Range("M91").Value = 404

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 297322.12
$ws.Range("I105").Value = 3096.5217
$ws.Range("K105").Value = 3096.5217
$ws.Range("M105").Value = -1349.5217

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 49999.332
$ws.Range("J23").Value = 49999
$ws.Range("L23").Value = 49999
$ws.Range("N23").Value = -50479

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 49999.332
$ws.Range("J27").Value = 49999
$ws.Range("L27").Value = 49999
$ws.Range("N27").Value = -50383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4113.6035
$ws.Range("I31").Value = 1333.7715
$ws.Range("J31").Value = 8343.781999999999
$ws.Range("K31").Value = 1333.7715
$ws.Range("L31").Value = 8343.781999999999
$ws.Range("M31").Value = -1038.7715
$ws.Range("N31").Value = -8933.781999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4113.6035
$ws.Range("I34").Value = 1333.7715
$ws.Range("J34").Value = 8343.781999999999
$ws.Range("K34").Value = 1333.7715
$ws.Range("L34").Value = 8343.781999999999
$ws.Range("M34").Value = -1131.7715
$ws.Range("N34").Value = -8747.781999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3063.3704
$ws.Range("I58").Value = 1234
$ws.Range("J58").Value = 4321.0625
$ws.Range("K58").Value = 1234
$ws.Range("L58").Value = 4321.0625
$ws.Range("M58").Value = -1031
$ws.Range("N58").Value = -4727.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 23454.6
$ws.Range("I62").Value = 42380
$ws.Range("J62").Value = 4529.2
$ws.Range("K62").Value = 42380
$ws.Range("L62").Value = 4529.2
$ws.Range("M62").Value = -41756
$ws.Range("N62").Value = -5777.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 23454.6
$ws.Range("I65").Value = 42380
$ws.Range("J65").Value = 4529.2
$ws.Range("K65").Value = 211900
$ws.Range("L65").Value = 22646
$ws.Range("M65").Value = -208780
$ws.Range("N65").Value = -28886

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 27000
$ws.Range("J70").Value = 27000
$ws.Range("L70").Value = 27000
$ws.Range("N70").Value = -27630

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 27000
$ws.Range("J73").Value = 27000
$ws.Range("L73").Value = 27000
$ws.Range("N73").Value = -29184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4888.7036
$ws.Range("I99").Value = 6513
$ws.Range("J99").Value = 2858.3333
$ws.Range("K99").Value = 6513
$ws.Range("L99").Value = 2858.3333
$ws.Range("M99").Value = -5015
$ws.Range("N99").Value = -5854.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4888.7036
$ws.Range("I126").Value = 6513
$ws.Range("J126").Value = 2858.3333
$ws.Range("K126").Value = 19539
$ws.Range("L126").Value = 8574.999899999999
$ws.Range("M126").Value = -17069
$ws.Range("N126").Value = -13514.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3063.3704
$ws.Range("I136").Value = 1234
$ws.Range("J136").Value = 4321.0625
$ws.Range("K136").Value = 3702
$ws.Range("L136").Value = 12963.1875
$ws.Range("M136").Value = -1152
$ws.Range("N136").Value = -18063.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1224.619
$ws.Range("I122").Value = 299.8
$ws.Range("J122").Value = 2065.3635
$ws.Range("K122").Value = 2698.2
$ws.Range("L122").Value = 18588.2715
$ws.Range("M122").Value = -248.2000000000003
$ws.Range("N122").Value = -23488.2715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6804141
$ws.Range("J131").Value = 7938072.5
$ws.Range("L131").Value = 23814217.5
$ws.Range("N131").Value = -23824297.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6315870
$ws.Range("I137").Value = 12502178
$ws.Range("J137").Value = 129562.375
$ws.Range("K137").Value = 37506534
$ws.Range("L137").Value = 388687.125
$ws.Range("M137").Value = -37501434
$ws.Range("N137").Value = -398887.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 336500
$ws.Range("I15").Value = 1000000
$ws.Range("J15").Value = 4750
$ws.Range("K15").Value = 1000000
$ws.Range("L15").Value = 4750
$ws.Range("M15").Value = -999712
$ws.Range("N15").Value = -5326

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 31700.2
$ws.Range("J39").Value = 31700.2
$ws.Range("L39").Value = 31700.2
$ws.Range("N39").Value = -32764.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6450.826
$ws.Range("I70").Value = 6631.6113
$ws.Range("J70").Value = 5800
$ws.Range("K70").Value = 6631.6113
$ws.Range("L70").Value = 5800
$ws.Range("M70").Value = -6361.6113
$ws.Range("N70").Value = -6340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6450.826
$ws.Range("I73").Value = 6631.6113
$ws.Range("J73").Value = 5800
$ws.Range("K73").Value = 6631.6113
$ws.Range("L73").Value = 5800
$ws.Range("M73").Value = -5695.6113
$ws.Range("N73").Value = -7672

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2860
$ws.Range("I80").Value = 2831.5625
$ws.Range("J80").Value = 2973.75
$ws.Range("K80").Value = 2831.5625
$ws.Range("L80").Value = 2973.75
$ws.Range("M80").Value = -1833.5625
$ws.Range("N80").Value = -4969.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 336500
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 4750
$ws.Range("K81").Value = 1000000
$ws.Range("L81").Value = 4750
$ws.Range("M81").Value = -999002
$ws.Range("N81").Value = -6746

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2860
$ws.Range("I83").Value = 2831.5625
$ws.Range("J83").Value = 2973.75
$ws.Range("K83").Value = 14157.8125
$ws.Range("L83").Value = 14868.75
$ws.Range("M83").Value = -9165.8125
$ws.Range("N83").Value = -24852.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 336500
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 4750
$ws.Range("K84").Value = 3000000
$ws.Range("L84").Value = 14250
$ws.Range("M84").Value = -2995008
$ws.Range("N84").Value = -24234

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 62149.668
$ws.Range("J124").Value = 62149.668
$ws.Range("L124").Value = 62149.668
$ws.Range("N124").Value = -71969.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 500118
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 500118
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3652.4
$ws.Range("I81").Value = 2205.5
$ws.Range("J81").Value = 4617
$ws.Range("K81").Value = 4411
$ws.Range("L81").Value = 9234
$ws.Range("M81").Value = -3350
$ws.Range("N81").Value = -11356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3652.4
$ws.Range("I84").Value = 2205.5
$ws.Range("J84").Value = 4617
$ws.Range("K84").Value = 22055
$ws.Range("L84").Value = 46170
$ws.Range("M84").Value = -16751
$ws.Range("N84").Value = -56778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 23233
$ws.Range("J101").Value = 23233
$ws.Range("L101").Value = 23233
$ws.Range("N101").Value = -29723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 63564.625
$ws.Range("I126").Value = 143419.86
$ws.Range("J126").Value = 1455
$ws.Range("K126").Value = 430259.58
$ws.Range("L126").Value = 4365
$ws.Range("M126").Value = -427789.58
$ws.Range("N126").Value = -9305

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 500000
$ws.Range("J135").Value = 500000
$ws.Range("L135").Value = 500000
$ws.Range("N135").Value = -510140
